$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet (LOGT2 -> logt2)
$ws.Name = "logt2"

# 2. Update the cycle "Stop" date in F1
$ws.Range("F1").Value = 41909

# 3. Insert a new row above the old header row (row 5), pushing the header
#    (old row 5) to row 6 and the existing data row (old row 6) to row 7.
$ws.Rows.Item(5).Insert()

# 4. Add the "Phase/Task" value for the existing (now shifted) data row.
$ws.Range("F7").Value = 12
$ws.Range("F7").HorizontalAlignment = -4152
$ws.Range("F7").Font.Name = "Times New Roman"
$ws.Range("F7").Font.Size = 11

# 5. New data row 8
$ws.Range("A8").Value = 41912
$ws.Range("B8").Value = 0.96875
$ws.Range("C8").Value = 0.993055555555556
$ws.Range("D8").Value = 0
$ws.Range("E8").Formula = "=((HOUR(C8)-HOUR(B8))*60)+(MINUTE(C8)-MINUTE(B8))-D8"
$ws.Range("F8").HorizontalAlignment = -4152
$ws.Range("F8").Font.Name = "Times New Roman"
$ws.Range("F8").Font.Size = 11
$ws.Range("H8").Value = "Tuvimos una reunión de equipo para la asignación de las tareas del ciclo #1."

# 6. New data row 9
$ws.Range("A9").Value = 41913
$ws.Range("B9").Value = 0.84375
$ws.Range("C9").Value = 0.927083333333333
$ws.Range("D9").Value = 0
$ws.Range("E9").Formula = "=((HOUR(C9)-HOUR(B9))*60)+(MINUTE(C9)-MINUTE(B9))-D9"
$ws.Range("F9").Value = "13,14,15"
$ws.Range("F9").HorizontalAlignment = -4152
$ws.Range("F9").Font.Name = "Times New Roman"
$ws.Range("F9").Font.Size = 11
$ws.Range("H9").Value = "Participe en el diagrama de use case y escenario de atributo de calidad"

# 7. New data row 10
$ws.Range("A10").Value = 41914
$ws.Range("B10").Value = 0.34375
$ws.Range("C10").Value = 0.427083333333333
$ws.Range("D10").Value = 45
$ws.Range("E10").Formula = "=((HOUR(C10)-HOUR(B10))*60)+(MINUTE(C10)-MINUTE(B10))-D10"
$ws.Range("F10").HorizontalAlignment = -4152
$ws.Range("F10").Font.Name = "Times New Roman"
$ws.Range("F10").Font.Size = 11
$ws.Range("H10").Value = "Documentar los inconvenientes encontrados"

# 8. Total row formula now lives in the newly inserted row 5 and sums the
#    (now expanded) data block.
$ws.Range("E5").Formula = "=SUM(E7:E10)/60"

# 9. Match the row heights used by the source workbook for the blank total
#    row and the three new (wrapped-comment) data rows.
$ws.Rows.Item(5).RowHeight = 13.75
$ws.Rows.Item(8).RowHeight = 26.65
$ws.Rows.Item(9).RowHeight = 26.65
$ws.Rows.Item(10).RowHeight = 14.15

# 10. Restore the selected cell shown when the sheet is opened.
$ws.Range("E6").Select()
